$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 44425
$ws.Range("J2").Value = 120
$ws.Range("K2").Value = 14000
$ws.Range("L2").Value = 15000
$ws.Range("M2").Value = 14500
$ws.Range("O2").Value = "Región del Maule"
$ws.Range("P2").Value = 362

# Row 3
$ws.Range("D3").Value = 44446
$ws.Range("J3").Value = 160
$ws.Range("K3").Value = 12500
$ws.Range("L3").Value = 13000
$ws.Range("M3").Value = 12750
$ws.Range("O3").Value = "Provincia del Elquí"
$ws.Range("P3").Value = 319

# Row 4
$ws.Range("D4").Value = 44454
$ws.Range("J4").Value = 120
$ws.Range("K4").Value = 13000
$ws.Range("L4").Value = 14000
$ws.Range("M4").Value = 13500
$ws.Range("P4").Value = 338

# Row 5
$ws.Range("D5").Value = 44489
$ws.Range("J5").Value = 120

# Row 6
$ws.Range("D6").Value = 44503
$ws.Range("J6").Value = 160
$ws.Range("K6").Value = 11000
$ws.Range("L6").Value = 12000
$ws.Range("M6").Value = 11500
$ws.Range("P6").Value = 288

# Row 7
$ws.Range("D7").Value = 44417
$ws.Range("K7").Value = 15000
$ws.Range("L7").Value = 16000
$ws.Range("M7").Value = 15500
$ws.Range("P7").Value = 388

# Row 8
$ws.Range("D8").Value = 44399
$ws.Range("H8").Value = "Española"
$ws.Range("I8").Value = "Segunda"
$ws.Range("K8").Value = 15500
$ws.Range("L8").Value = 16000
$ws.Range("M8").Value = 15750
$ws.Range("P8").Value = 394

# Row 9
$ws.Range("D9").Value = 44475

# Row 10
$ws.Range("D10").Value = 44510
$ws.Range("K10").Value = 11000
$ws.Range("L10").Value = 12000
$ws.Range("M10").Value = 11500
$ws.Range("O10").Value = "Provincia del Elquí"
$ws.Range("P10").Value = 288

# Row 11
$ws.Range("D11").Value = 44512

# Row 12
$ws.Range("D12").Value = 44467
$ws.Range("J12").Value = 160
$ws.Range("O12").Value = "Provincia de Limarí"

# Row 13
$ws.Range("D13").Value = 44508
$ws.Range("J13").Value = 160
$ws.Range("K13").Value = 11000
$ws.Range("L13").Value = 12000
$ws.Range("M13").Value = 11500
$ws.Range("P13").Value = 288

# Row 14
$ws.Range("D14").Value = 44420
$ws.Range("J14").Value = 120
$ws.Range("K14").Value = 13000
$ws.Range("L14").Value = 14000
$ws.Range("M14").Value = 13500
$ws.Range("P14").Value = 338

# Row 15
$ws.Range("D15").Value = 44495

# Row 16
$ws.Range("D16").Value = 44498
$ws.Range("J16").Value = 60
$ws.Range("K16").Value = 10500
$ws.Range("L16").Value = 11000
$ws.Range("M16").Value = 10750
$ws.Range("P16").Value = 269

# Row 17
$ws.Range("D17").Value = 44427
$ws.Range("K17").Value = 13000
$ws.Range("L17").Value = 14000
$ws.Range("M17").Value = 13500
$ws.Range("P17").Value = 338

# Row 18
$ws.Range("D18").Value = 44468
$ws.Range("J18").Value = 60
$ws.Range("K18").Value = 12000
$ws.Range("M18").Value = 12500
$ws.Range("P18").Value = 312

# Row 19
$ws.Range("D19").Value = 44515
$ws.Range("J19").Value = 120
$ws.Range("K19").Value = 11000
$ws.Range("L19").Value = 12000
$ws.Range("M19").Value = 11500
$ws.Range("P19").Value = 288

# Row 20
$ws.Range("D20").Value = 44453
$ws.Range("J20").Value = 160
$ws.Range("K20").Value = 12500
$ws.Range("L20").Value = 13000
$ws.Range("M20").Value = 12750
$ws.Range("P20").Value = 319

# Row 21
$ws.Range("D21").Value = 44487
$ws.Range("J21").Value = 100
$ws.Range("K21").Value = 11000
$ws.Range("L21").Value = 12000
$ws.Range("M21").Value = 11500
$ws.Range("P21").Value = 288

# Row 22
$ws.Range("D22").Value = 44496
$ws.Range("H22").Value = "Madrigal"
$ws.Range("I22").Value = "Primera"
$ws.Range("K22").Value = 11000
$ws.Range("L22").Value = 12000
$ws.Range("M22").Value = 11500
$ws.Range("P22").Value = 288

# Row 23
$ws.Range("D23").Value = 44490
$ws.Range("J23").Value = 100

# Row 24
$ws.Range("D24").Value = 44432
$ws.Range("K24").Value = 14000
$ws.Range("L24").Value = 15000
$ws.Range("M24").Value = 14500
$ws.Range("P24").Value = 362

# Row 25
$ws.Range("D25").Value = 44484

# Row 26
$ws.Range("D26").Value = 44488
$ws.Range("J26").Value = 100
$ws.Range("K26").Value = 11000
$ws.Range("L26").Value = 12000
$ws.Range("M26").Value = 11500
$ws.Range("P26").Value = 288

# Row 27
$ws.Range("D27").Value = 44455
$ws.Range("K27").Value = 13000
$ws.Range("L27").Value = 14000
$ws.Range("M27").Value = 13500
$ws.Range("P27").Value = 338

# Row 28
$ws.Range("D28").Value = 44516
$ws.Range("K28").Value = 11000
$ws.Range("L28").Value = 12000
$ws.Range("M28").Value = 11500
$ws.Range("P28").Value = 288

# Row 29
$ws.Range("D29").Value = 44491
$ws.Range("J29").Value = 100
$ws.Range("K29").Value = 11000
$ws.Range("L29").Value = 12000
$ws.Range("M29").Value = 11500
$ws.Range("P29").Value = 288

# Row 30
$ws.Range("D30").Value = 44435
$ws.Range("J30").Value = 120
$ws.Range("K30").Value = 14000
$ws.Range("L30").Value = 15000
$ws.Range("M30").Value = 14500
$ws.Range("O30").Value = "Provincia del Elquí"
$ws.Range("P30").Value = 362

# Row 31
$ws.Range("D31").Value = 44482
$ws.Range("J31").Value = 120

# Row 32
$ws.Range("D32").Value = 44494
$ws.Range("J32").Value = 120
$ws.Range("M32").Value = 11500
$ws.Range("P32").Value = 288

# Row 33
$ws.Range("D33").Value = 44426
$ws.Range("K33").Value = 13000
$ws.Range("L33").Value = 14000
$ws.Range("M33").Value = 13500
$ws.Range("O33").Value = "Región del Maule"
$ws.Range("P33").Value = 338

# Row 34
$ws.Range("D34").Value = 44473
$ws.Range("J34").Value = 160

# Row 35
$ws.Range("D35").Value = 44505
